$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows worth of cells in columns A:D only (shifting existing
# A:D data down by 2 rows) - the H:O summary table to the right must stay put.
$ws.Range("A11:D12").Insert(-4121) # xlShiftDown

# Populate the two newly inserted rows with the "i 2 verbali esterni" data.
$ws.Cells.Item(11, 1).Value = "Viktorija"
$ws.Cells.Item(11, 2).Value = "Documenti di progetto"
$ws.Cells.Item(11, 3).Value = 75
$ws.Cells.Item(11, 4).Value = 43497

$ws.Cells.Item(12, 1).Value = "Viktorija"
$ws.Cells.Item(12, 2).Value = "Documenti di progetto"
$ws.Cells.Item(12, 3).Value = 55
$ws.Cells.Item(12, 4).Value = 43503

# Restore the view the workbook was saved with.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("E20").Select()
